$d = $word.ActiveDocument

# The "25/11/2012: GUI verbeterd..." entry is the last list item in the
# Logboek and carries the _GoBack bookmark at the end of its paragraph.
# Find its paragraph index, then insert a brand-new list paragraph (same
# ListParagraph style / numbering) right after it -- the bookmark stays
# where it is, at the end of the original paragraph -- holding the new
# log entry text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*25/11/2012: GUI verbeterd*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs($targetIndex)
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($targetIndex + 1)
$newPara.Range.Text = "25/11/2012: Geëxperimenteerd met json/xml. Men kan momenteel nu input geven en een echt! artist object wordt gereturned van de last.fm api. Ik extraheer en toon momenteel enkel de bio summary, maar moet nog een gui met scroll + foto + rating pagina/activity maken. "
